$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell 2 4 '42.076.90'
Set-TextCell 2 5 '  -1.08%  '
Set-TextCell 3 4 '2.260.27'
Set-TextCell 3 5 '  -1.14%  '
Set-TextCell 4 5 '  -0.16%  '
Set-TextCell 5 4 '306.89'
Set-TextCell 5 5 '  +0.14%  '
Set-TextCell 6 4 '96.86'
Set-TextCell 6 5 '  +0.26%  '
Set-TextCell 7 5 '  -1.62%  '
Set-TextCell 8 5 '  -0.07%  '
Set-TextCell 9 5 '  -1.65%  '
Set-TextCell 10 4 '35.06'
Set-TextCell 10 5 '  -3.72%  '
Set-TextCell 11 4 '0.0786'
Set-TextCell 11 5 '  -2.29%  '
Set-TextCell 12 5 '  +0.54%  '
Set-TextCell 13 4 '6.81'
Set-TextCell 13 5 '  +1.01%  '
Set-TextCell 14 4 '2.610.07'
Set-TextCell 14 5 '  -1.30%  '
Set-TextCell 15 5 '  -0.01%  '
Set-TextCell 16 4 '2.257.24'
Set-TextCell 16 5 '  -2.33%  '
Set-TextCell 17 4 '0.787'
Set-TextCell 17 5 '  -2.04%  '
Set-TextCell 18 4 '41.890.02'
Set-TextCell 18 5 '  -1.30%  '
Set-TextCell 19 4 '12.19'
Set-TextCell 20 4 '0.0₃0900'
Set-TextCell 20 5 '  -2.17%  '
Set-TextCell 21 4 '5.97'
Set-TextCell 21 5 '  -0.62%  '
Set-TextCell 22 4 '67.61'
Set-TextCell 22 5 '  -0.43%  '
Set-TextCell 23 5 '  -2.91%  '
Set-TextCell 24 4 '1.97'
Set-TextCell 24 5 '  +0.95%  '
Set-TextCell 25 4 '2.58'
Set-TextCell 25 5 '  -1.12%  '
Set-TextCell 26 5 '  -0.13%  '
Set-TextCell 27 4 '23.47'
Set-TextCell 27 5 '  -2.43%  '
Set-TextCell 28 4 '36.61'
Set-TextCell 28 5 '  +0.29%  '
Set-TextCell 29 5 '  +0.98%  '
Set-TextCell 30 4 '9.51'
Set-TextCell 30 5 '  -1.08%  '
Set-TextCell 31 4 '165.05'
Set-TextCell 31 5 '  +1.97%  '
Set-TextCell 32 4 '5.22'
Set-TextCell 32 5 '  -2.13%  '
Set-TextCell 33 5 '  -0.03%  '
Set-TextCell 34 4 '3.13'
Set-TextCell 34 5 '  +1.05%  '
Set-TextCell 35 4 '0.0735'
Set-TextCell 35 5 '  -2.38%  '
Set-TextCell 36 4 '17.43'
Set-TextCell 36 5 '  +0.05%  '
Set-TextCell 37 5 '  +0.45%  '
Set-TextCell 38 5 '  -4.49%  '
Set-TextCell 39 5 '  -1.28%  '
Set-TextCell 40 5 '  -3.70%  '
Set-TextCell 41 4 '4.12'
Set-TextCell 41 5 '  -2.15%  '
Set-TextCell 42 5 '  -2.38%  '
Set-TextCell 43 4 '1.949.22'
Set-TextCell 43 5 '  -2.81%  '
Set-TextCell 44 4 '18.95'
Set-TextCell 44 5 '  -2.12%  '
Set-TextCell 45 5 '  -1.93%  '
Set-TextCell 46 4 '2.91'
Set-TextCell 46 5 '  -3.44%  '
Set-TextCell 47 4 '9.80'
Set-TextCell 47 5 '  -4.18%  '
Set-TextCell 48 4 '53.39'
Set-TextCell 48 5 '  -1.39%  '
Set-TextCell 49 4 '2.482.37'
Set-TextCell 49 5 '  -1.21%  '
Set-TextCell 50 4 '92.45'
Set-TextCell 50 5 '  +0.40%  '
Set-TextCell 51 4 '71.52'
Set-TextCell 51 5 '  -1.74%  '
